$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monitorees")

# ---------------------------------------------------------------------------
# SARAALERT-1260: Allow vaccine table to be populated on import
# Adds two vaccine blocks (Group Name, Product Name, Administration Date,
# Dose Number, Notes) as new trailing columns CY:DH on the Monitorees sheet.
# ---------------------------------------------------------------------------

# Set text number format on the Administration Date columns (DA, DF) BEFORE writing values
# so the date-like strings are stored as text, matching columns 105/110 style in the target.
# Only the rows that actually receive vaccine data get the "@" text format.
$ws.Range("DA1:DA7").NumberFormat = "@"
$ws.Range("DF1:DF4").NumberFormat = "@"

# Row 1
$ws.Range("CY1").Value = "Vaccine 1 Group Name"
$ws.Range("CZ1").Value = "Vaccine 1 Product Name"
$ws.Range("DA1").Value = "Vaccine 1 Administration Date"
$ws.Range("DB1").Value = "Vaccine 1 Dose Number"
$ws.Range("DC1").Value = "Vaccine 1 Notes"
$ws.Range("DD1").Value = "Vaccine 2 Group Name"
$ws.Range("DE1").Value = "Vaccine 2 Product Name"
$ws.Range("DF1").Value = "Vaccine 2 Administration Date"
$ws.Range("DG1").Value = "Vaccine 2 Dose Number"
$ws.Range("DH1").Value = "Vaccine 2 Notes"

# Row 2
$ws.Range("CY2").Value = "COVID-19"
$ws.Range("CZ2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DA2").Value = "2020-06-01"
$ws.Range("DB2").Value = 1
$ws.Range("DC2").Value = "notes 1"
$ws.Range("DD2").Value = "COVID-19"
$ws.Range("DE2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DF2").Value = "2020-06-20"
$ws.Range("DG2").Value = 2
$ws.Range("DH2").Value = "notes 2"

# Row 3
$ws.Range("CY3").Value = "COVID-19"
$ws.Range("CZ3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("DA3").Value = "2020-06-02"
$ws.Range("DB3").Value = 1
$ws.Range("DD3").Value = "COVID-19"
$ws.Range("DE3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("DF3").Value = "2020-06-21"
$ws.Range("DG3").Value = 2

# Row 4
$ws.Range("CY4").Value = "COVID-19"
$ws.Range("CZ4").Value = "Unknown"
$ws.Range("DA4").Value = "2020-06-04"
$ws.Range("DB4").Value = 1
$ws.Range("DD4").Value = "COVID-19"
$ws.Range("DE4").Value = "Unknown"
$ws.Range("DF4").Value = "2020-06-22"
$ws.Range("DG4").Value = 2

# Row 5
$ws.Range("CY5").Value = "COVID-19"
$ws.Range("CZ5").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DA5").Value = "2020-06-01"
$ws.Range("DB5").Value = 1

# Row 6
$ws.Range("CY6").Value = "COVID-19"
$ws.Range("CZ6").Value = "Janssen (J&J) COVID-19 Vaccine"
$ws.Range("DA6").Value = "2020-06-03"
$ws.Range("DB6").Value = 1

# Row 7
$ws.Range("CY7").Value = "COVID-19"
$ws.Range("CZ7").Value = "Unknown"
$ws.Range("DA7").Value = "2020-06-02"
$ws.Range("DB7").Value = 1
# ---------------------------------------------------------------------------
# Column widths for the newly added vaccine columns (CY:DH), matching the
# widths used by the other "bestFit" columns on this sheet.
# ---------------------------------------------------------------------------
$ws.Columns.Item(103).ColumnWidth = 20.33203125  # CY  Vaccine 1 Group Name
$ws.Columns.Item(104).ColumnWidth = 31            # CZ  Vaccine 1 Product Name
$ws.Columns.Item(105).ColumnWidth = 25.6640625   # DA  Vaccine 1 Administration Date
$ws.Columns.Item(106).ColumnWidth = 21.1640625   # DB  Vaccine 1 Dose Number
$ws.Columns.Item(107).ColumnWidth = 14.5          # DC  Vaccine 1 Notes
$ws.Columns.Item(108).ColumnWidth = 20.33203125  # DD  Vaccine 2 Group Name
$ws.Columns.Item(109).ColumnWidth = 31            # DE  Vaccine 2 Product Name
$ws.Columns.Item(110).ColumnWidth = 25.6640625   # DF  Vaccine 2 Administration Date
$ws.Columns.Item(111).ColumnWidth = 21.1640625   # DG  Vaccine 2 Dose Number
$ws.Columns.Item(112).ColumnWidth = 14.5          # DH  Vaccine 2 Notes

# ---------------------------------------------------------------------------
# Reset the sheet view: scroll back to the top-left (A1) and select A1,
# instead of the prior scrolled view (topLeftCell=CM1, selection=CY9).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A1").Select()
